# Auto-generated edit script applying the Ridill_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 15207.8
$ws.Range("I53").Value = 37557.5
$ws.Range("J53").Value = 308
$ws.Range("K53").Value = 37557.5
$ws.Range("L53").Value = 308
$ws.Range("M53").Value = -36920.5
$ws.Range("N53").Value = -1582
$ws.Range("H82").Value = 5322
$ws.Range("J82").Value = 7600
$ws.Range("L82").Value = 22800
$ws.Range("N82").Value = -23612
$ws.Range("H85").Value = 5322
$ws.Range("J85").Value = 7600
$ws.Range("L85").Value = 22800
$ws.Range("N85").Value = -25608
$ws.Range("H103").Value = 45456812
$ws.Range("I103").Value = 100001200
$ws.Range("J103").Value = 3158.75
$ws.Range("K103").Value = 300003600
$ws.Range("L103").Value = 9476.25
$ws.Range("M103").Value = -300003014
$ws.Range("N103").Value = -10648.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 455789.97
$ws.Range("J45").Value = 1029
$ws.Range("L45").Value = 1029
$ws.Range("N45").Value = -1783
$ws.Range("H61").Value = 5583675
$ws.Range("I61").Value = 2605101.5
$ws.Range("J61").Value = 29412264
$ws.Range("K61").Value = 2605101.5
$ws.Range("L61").Value = 29412264
$ws.Range("M61").Value = -2604889.5
$ws.Range("N61").Value = -29412688
$ws.Range("H74").Value = 67290376
$ws.Range("I74").Value = 72450070
$ws.Range("J74").Value = 59264176
$ws.Range("K74").Value = 72450070
$ws.Range("L74").Value = 59264176
$ws.Range("M74").Value = -72449196
$ws.Range("N74").Value = -59265924
$ws.Range("H77").Value = 67290376
$ws.Range("I77").Value = 72450070
$ws.Range("J77").Value = 59264176
$ws.Range("K77").Value = 362250350
$ws.Range("L77").Value = 296320880
$ws.Range("M77").Value = -362245982
$ws.Range("N77").Value = -296329616
$ws.Range("H94").Value = 33333
$ws.Range("J94").Value = 33333
$ws.Range("L94").Value = 33333
$ws.Range("N94").Value = -35135
$ws.Range("H106").Value = 38000
$ws.Range("J106").Value = 38000
$ws.Range("L106").Value = 38000
$ws.Range("N106").Value = -40524
$ws.Range("H136").Value = 5583675
$ws.Range("I136").Value = 2605101.5
$ws.Range("J136").Value = 29412264
$ws.Range("K136").Value = 7815304.5
$ws.Range("L136").Value = 88236792
$ws.Range("M136").Value = -7812754.5
$ws.Range("N136").Value = -88241892

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1904.26
$ws.Range("I86").Value = 1930.1158
$ws.Range("J86").Value = 1413
$ws.Range("K86").Value = 1930.1158
$ws.Range("L86").Value = 1413
$ws.Range("M86").Value = -807.1158
$ws.Range("N86").Value = -3659
$ws.Range("H89").Value = 1904.26
$ws.Range("I89").Value = 1930.1158
$ws.Range("J89").Value = 1413
$ws.Range("K89").Value = 9650.579
$ws.Range("L89").Value = 7065
$ws.Range("M89").Value = -4034.579
$ws.Range("N89").Value = -18297
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H105").Value = 2212.6667
$ws.Range("I105").Value = 2112.8572
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 2112.8572
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = -365.8571999999999
$ws.Range("N105").Value = -5794
$ws.Range("H134").Value = 19843206
$ws.Range("I134").Value = 26317260
$ws.Range("J134").Value = 4467325
$ws.Range("K134").Value = 78951780
$ws.Range("L134").Value = 13401975
$ws.Range("M134").Value = -78949245
$ws.Range("N134").Value = -13407045

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2621355.8
$ws.Range("I58").Value = 1556139
$ws.Range("J58").Value = 5683854
$ws.Range("K58").Value = 1556139
$ws.Range("L58").Value = 5683854
$ws.Range("M58").Value = -1555936
$ws.Range("N58").Value = -5684260
$ws.Range("H136").Value = 2621355.8
$ws.Range("I136").Value = 1556139
$ws.Range("J136").Value = 5683854
$ws.Range("K136").Value = 4668417
$ws.Range("L136").Value = 17051562
$ws.Range("M136").Value = -4665867
$ws.Range("N136").Value = -17056662

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1765.4
$ws.Range("I17").Value = 262.5
$ws.Range("J17").Value = 7777
$ws.Range("K17").Value = 787.5
$ws.Range("L17").Value = 23331
$ws.Range("M17").Value = -618.5
$ws.Range("N17").Value = -23669
$ws.Range("H131").Value = 959.6316
$ws.Range("I131").Value = 435.85715
$ws.Range("J131").Value = 1265.1666
$ws.Range("K131").Value = 1307.57145
$ws.Range("L131").Value = 3795.4998
$ws.Range("M131").Value = 3732.42855
$ws.Range("N131").Value = -13875.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6017.8945
$ws.Range("I102").Value = 8546.666999999999
$ws.Range("J102").Value = 1682.8572
$ws.Range("K102").Value = 8546.666999999999
$ws.Range("L102").Value = 1682.8572
$ws.Range("M102").Value = -6924.666999999999
$ws.Range("N102").Value = -4926.8572
$ws.Range("H132").Value = 15390357
$ws.Range("I132").Value = 16509652
$ws.Range("K132").Value = 49528956
$ws.Range("M132").Value = -49526426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1290
$ws.Range("I7").Value = 780
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 780
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -668
$ws.Range("N7").Value = -2024
$ws.Range("H40").Value = 2982.5557
$ws.Range("I40").Value = 3200
$ws.Range("J40").Value = 2920.4285
$ws.Range("K40").Value = 3200
$ws.Range("L40").Value = 2920.4285
$ws.Range("M40").Value = -3064
$ws.Range("N40").Value = -3192.4285
$ws.Range("H122").Value = 18441864
$ws.Range("I122").Value = 1936579.6
$ws.Range("J122").Value = 200000000
$ws.Range("K122").Value = 5809738.800000001
$ws.Range("L122").Value = 600000000
$ws.Range("M122").Value = -5807288.800000001
$ws.Range("N122").Value = -600004900
$ws.Range("H126").Value = 1290
$ws.Range("I126").Value = 780
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 2340
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = 130
$ws.Range("N126").Value = -10340
$ws.Range("H132").Value = 16670292
$ws.Range("I132").Value = 66666664
$ws.Range("J132").Value = 4834.3335
$ws.Range("K132").Value = 199999992
$ws.Range("L132").Value = 14503.0005
$ws.Range("M132").Value = -199997462
$ws.Range("N132").Value = -19563.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 41672370
$ws.Range("I126").Value = 250000000
$ws.Range("J126").Value = 6839.8
$ws.Range("K126").Value = 750000000
$ws.Range("L126").Value = 20519.4
$ws.Range("M126").Value = -749997530
$ws.Range("N126").Value = -25459.4
$ws.Range("H132").Value = 3088939.8
$ws.Range("I132").Value = 2303275.2
$ws.Range("K132").Value = 6909825.600000001
$ws.Range("M132").Value = -6907295.600000001
$ws.Range("H136").Value = 18249.77
$ws.Range("I136").Value = 11567
$ws.Range("J136").Value = 55005
$ws.Range("K136").Value = 34701
$ws.Range("L136").Value = 165015
$ws.Range("M136").Value = -32151
$ws.Range("N136").Value = -170115
